# Update 想去人数 (F) / 最低票价 (G) figures on the "展览" and "全部类型"
# sheets to match the latest scrape (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")

# Map of row => @{ F = value; G = value } for each worksheet that needs updating.
$exhibitUpdates = @(
    @{ Row = 2;  F = 1930; G = 70 },
    @{ Row = 3;  G = 65 },
    @{ Row = 4;  F = 172 },
    @{ Row = 5;  F = 185 },
    @{ Row = 6;  F = 2872 },
    @{ Row = 7;  F = 201 },
    @{ Row = 8;  F = 101 },
    @{ Row = 10; F = 1598 },
    @{ Row = 11; F = 647 },
    @{ Row = 23; F = 27 },
    @{ Row = 24; F = 1 },
    @{ Row = 25; F = 271 },
    @{ Row = 26; F = 85 },
    @{ Row = 27; F = 86 },
    @{ Row = 29; F = 1867 },
    @{ Row = 31; F = 437 },
    @{ Row = 33; F = 121 },
    @{ Row = 38; F = 467 }
)

$allUpdates = @(
    @{ Row = 2;  F = 1930; G = 70 },
    @{ Row = 4;  G = 65 },
    @{ Row = 5;  F = 172 },
    @{ Row = 6;  F = 185 },
    @{ Row = 7;  F = 2872 },
    @{ Row = 8;  F = 201 },
    @{ Row = 9;  F = 101 },
    @{ Row = 11; F = 1598 },
    @{ Row = 12; F = 648 },
    @{ Row = 24; F = 27 },
    @{ Row = 25; F = 1 },
    @{ Row = 26; F = 271 },
    @{ Row = 27; F = 85 },
    @{ Row = 28; F = 86 },
    @{ Row = 30; F = 1867 },
    @{ Row = 32; F = 437 },
    @{ Row = 34; F = 121 },
    @{ Row = 39; F = 467 }
)

foreach ($u in $exhibitUpdates) {
    if ($u.ContainsKey('F')) {
        $wsExhibit.Range("F$($u.Row)").Value = $u.F
    }
    if ($u.ContainsKey('G')) {
        $wsExhibit.Range("G$($u.Row)").Value = $u.G
    }
}

foreach ($u in $allUpdates) {
    if ($u.ContainsKey('F')) {
        $wsAll.Range("F$($u.Row)").Value = $u.F
    }
    if ($u.ContainsKey('G')) {
        $wsAll.Range("G$($u.Row)").Value = $u.G
    }
}
